$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 12 (shifts existing rows 12..30 down to 13..31).
$ws.Rows.Item(12).Insert()

# Fill in the newly inserted row: new "Programming Project 4" due 2024-02-16 (45338).
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Programming Project 4"
$ws.Range("C12").Value = 45338

# Renumber every "Programming Project N" label from row 13 down to the last
# data row (now row 30) up by one, since the new project took the "4" slot.
for ($r = 13; $r -le 30; $r++) {
    $label = $ws.Range("B$r").Value()
    if ($label -match '^Programming Project (\d+)$') {
        $n = [int]$matches[1]
        $ws.Range("B$r").Value = "Programming Project $($n + 1)"
    }
}

# Re-sequence column A (the "number" column) for every data row so it stays 1..29.
for ($r = 2; $r -le 30; $r++) {
    $ws.Range("A$r").Value = $r - 1
}

# Keep the selection / active cell in sync with the now-longer sheet, matching
# the shifted trailing blank row.
$ws.Range("A31").Select()
